$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '61.764.97'
$ws.Range('D2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.407.13'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '412.84'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '129.11'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('E7').Value = '  -2.87%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('E10').Value = '  -4.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.65'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.949.19'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '12.70'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +4.57%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.398.28'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '61.801.76'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '478.98'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +10.49%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '90.60'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.27'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.48%  '
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.31'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.79'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +10.85%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '33.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '4.76'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.71'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.07%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '11.84'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.64'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.31%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.166'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.111'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '40.84'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '58.35'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +8.03%  '
$ws.Range('E37').Value = '  -2.41%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  +4.10%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '149.03'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.74%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.326'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +3.68%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('E44').Value = '  +4.90%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.24'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.58%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.57'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +6.44%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.36'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +19.51%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '16.39'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.91%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0₃0544'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +22.55%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '22.21'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '113.22'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +10.64%  '
